$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "  -3.64%  "
    3  = "  -5.43%  "
    4  = "  +0.07%  "
    5  = "  -5.48%  "
    6  = "  -6.59%  "
    7  = "  +0.00%  "
    8  = "  -5.34%  "
    9  = "  -6.22%  "
    10 = "  -7.44%  "
    11 = "  -10.44%  "
    12 = "  -9.51%  "
    13 = "  -4.90%  "
    14 = "  -1.22%  "
    15 = "  -6.94%  "
    16 = "  -6.21%  "
    17 = "  -5.52%  "
    18 = "  -9.48%  "
    19 = "  -6.88%  "
    20 = "  -6.50%  "
    21 = "  -7.67%  "
    22 = "  -9.82%  "
    23 = "  -0.21%  "
    24 = "  -7.15%  "
    25 = "  -5.99%  "
    26 = "  -2.53%  "
    27 = "  -0.59%  "
    28 = "  -8.67%  "
    29 = "  +0.03%  "
    30 = "  -5.53%  "
    31 = "  +0.08%  "
    32 = "  -7.82%  "
    33 = "  -5.89%  "
    34 = "  -6.38%  "
    35 = "  -5.01%  "
    36 = "  -6.66%  "
    37 = "  -8.48%  "
    38 = "  -10.99%  "
    39 = "  -6.14%  "
    40 = "  -6.75%  "
    41 = "  -6.04%  "
    42 = "  -11.82%  "
    43 = "  +0.08%  "
    44 = "  -8.55%  "
    45 = "  -4.25%  "
    46 = "  -7.48%  "
    47 = "  -2.33%  "
    48 = "  -8.99%  "
    49 = "  -5.24%  "
    50 = "  -6.84%  "
    51 = "  -7.68%  "
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
